$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Repull data / push all data / mean calculation:
# Updated "dSF" (column F) values for several rows to reflect repulled data.
$ws.Range("F3").Value = 5
$ws.Range("F5").Value = 9
$ws.Range("F6").Value = -3
$ws.Range("F8").Value = 8
$ws.Range("F9").Value = 5
$ws.Range("F12").Value = -2
$ws.Range("F16").Value = 0
$ws.Range("F18").Value = 10
$ws.Range("F22").Value = -6
$ws.Range("F25").Value = 7
